$d = $word.ActiveDocument

$replacements = @(
    @("58÷7=", "11÷9="),
    @("14÷9=", "16÷4="),
    @("82÷7=", "38÷5="),
    @("20÷3=", "50÷9="),
    @("87÷9=", "24÷2="),
    @("33÷5=", "17÷7="),
    @("14÷2=", "89÷6="),
    @("81÷4=", "33÷9="),
    @("69÷5=", "54÷2="),
    @("30÷9=", "11÷4="),
    @("90÷6=", "28÷7="),
    @("80÷7=", "25÷6="),
    @("86÷2=", "97÷7="),
    @("43÷4=", "72÷4="),
    @("79÷9=", "14÷7="),
    @("71÷9=", "99÷7="),
    @("19÷3=", "11÷8="),
    @("81÷7=", "16÷8="),
    @("71÷5=", "66÷3="),
    @("95÷9=", "86÷5="),
    @("55÷7=", "59÷2="),
    @("77÷7=", "39÷7="),
    @("89÷5=", "78÷5="),
    @("73÷4=", "96÷7="),
    @("57÷2=", "95÷3=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
